# Consolidate the run-per-word text runs in the "Followed by a picture"
# caption textbox into fewer runs, each carrying a trailing space, e.g.
#   "Followed" + " " + "by" + " " + "a" + " " + "picture"
# becomes
#   "Followed " + "by " + "a " + "picture"
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Merge "Followed" + " " (chars 1-9) into a single run "Followed ".
$tr.Characters(1, 9).Text = "Followed "

# Merge "by" + " " (chars 10-12) into a single run "by ".
$tr.Characters(10, 3).Text = "by "

# Merge "a" + " " (chars 13-14) into a single run "a ".
$tr.Characters(13, 2).Text = "a "

# The trailing "picture" run is untouched.
